$d = $word.ActiveDocument

# 1) Remove the accidental duplicated wording:
#      " project more that was more than just"  ->  " project more than just"
#    ("project more tha" + "t was more than just" collapses to "project more than just")
$found1 = $d.Content.Find.Execute(
    " project more that was more than just", # FindText
    $true,                                    # MatchCase
    $false,                                   # MatchWholeWord
    $false,                                   # MatchWildcards
    $false,                                   # MatchSoundsLike
    $false,                                   # MatchAllWordForms
    $true,                                    # Forward
    1,                                         # Wrap (wdFindContinue)
    $false,                                   # Format
    " project more than just",                # ReplaceWith
    2                                          # Replace (wdReplaceAll)
)
if (-not $found1) {
    throw "edit.ps1: could not find/replace the 'project more that was more than just' text"
}

# 2) Reword the sentence that introduces the technique list:
#      " To address the overall topic, the following techniques were applied:"
#      -> " The following techniques were applied to address the topic:"
$found2 = $d.Content.Find.Execute(
    " To address the overall topic, the following techniques were applied:", # FindText
    $true,                                    # MatchCase
    $false,                                   # MatchWholeWord
    $false,                                   # MatchWildcards
    $false,                                   # MatchSoundsLike
    $false,                                   # MatchAllWordForms
    $true,                                    # Forward
    1,                                         # Wrap (wdFindContinue)
    $false,                                   # Format
    " The following techniques were applied to address the topic:", # ReplaceWith
    2                                          # Replace (wdReplaceAll)
)
if (-not $found2) {
    throw "edit.ps1: could not find/replace the 'To address the overall topic...' text"
}
